$d = $word.ActiveDocument

# Locate the target list item by its text using Find.
$searchText = "Write out array of objects for each step/choice during the gameplay (components will reference the text values at specific indices to know what needs to be displayed at what time)"

$rng = $d.Content
$found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Target paragraph not found"
}

# Expand the range to the whole paragraph (including the paragraph mark)
# so that both the run's rPr and the paragraph mark's rPr (w:pPr/w:rPr)
# receive the strikethrough formatting - matching the behavior of
# selecting the complete paragraph in Word and applying Strikethrough.
$para = $rng.Paragraphs(1)
$paraRange = $para.Range
$paraRange.Font.Strikethrough = $true
